# Delete the KODK row (row 195, ticker KODK, fixed_quarter_date 2016-10-01)
# from the worksheet. All rows below shift up by one, and the used range
# shrinks from A1:D403 to A1:D402.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(195).Delete()
